# Edit script for baddata2015test.xlsx
# Summary of changes:
#  - Insert a new column G ("birthdate") before the existing bib_number column,
#    shifting bib_number..Finish from G:O to H:P.
#  - Populate the new birthdate column for the first two data rows with dates,
#    formatted as date (d-mmm-yy / builtin numFmt 15).
#  - Correct several "age" (column F) values.
#  - Expand the print area from A1:J6 to A1:K6.
#  - Move the active selection to G8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at G (pushes bib_number..Finish right by one).
$ws.Columns("G:G").Insert()

# Header for the freshly inserted column.
$ws.Cells.Item(1, 7).Value = "birthdate"

# Birthdates for the first two runners (stored as date serials).
$ws.Range("G3").Value = 31809
$ws.Range("G3").NumberFormat = "d-mmm-yy"

$ws.Range("G4").Value = 30390
$ws.Range("G4").NumberFormat = "d-mmm-yy"

# Fix up several age values in column F.
$ws.Range("F6").Value = 25
$ws.Range("F7").Value = 51
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = ""

# Widen the print area now that there is an extra column of data.
$ws.PageSetup.PrintArea = '$A$1:$K$6'

# Leave the selection where the editor last left it.
$ws.Range("G8").Select()
